$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 34.67150079277169
$ws.Range("G2").Value = 32.308185612594
$ws.Range("H2").Value = 37.08203046344973
$ws.Range("I2").Value = 0.004875439521266707
$ws.Range("J2").Value = 0.0007247618761747751
$ws.Range("K2").Value = 0.01314330250379633
$ws.Range("L2").Value = 0.0569748011382435
$ws.Range("M2").Value = 0.05463332696829905
$ws.Range("N2").Value = 0.05991939919198197

# Row 3
$ws.Range("F3").Value = 0.5265381850453629
$ws.Range("G3").Value = 0.02048527309711148
$ws.Range("H3").Value = 1.038044166319508
$ws.Range("I3").Value = 0.4916629681751183
$ws.Range("J3").Value = 0.0192122726768235
$ws.Range("K3").Value = 0.9691492111444576
$ws.Range("L3").Value = 0.5465883663268607
$ws.Range("M3").Value = 0.02119892645419704
$ws.Range("N3").Value = 1.076694064868029

# Row 4
$ws.Range("F4").Value = 35.19803897781706
$ws.Range("G4").Value = 32.32867088569111
$ws.Range("H4").Value = 38.12007462976923
$ws.Range("I4").Value = 0.496538407696385
$ws.Range("J4").Value = 0.01993703455299828
$ws.Range("K4").Value = 0.9822925136482539
$ws.Range("L4").Value = 0.6035631674651043
$ws.Range("M4").Value = 0.07583225342249608
$ws.Range("N4").Value = 1.136613464060011
